$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 139
$ws.Range("I9").Value = 157.8
$ws.Range("J9").Value = 120.2
$ws.Range("K9").Value = 157.8
$ws.Range("L9").Value = 120.2
$ws.Range("M9").Value = 11.19999999999999
$ws.Range("N9").Value = -458.2

$ws.Range("H28").Value = 198.74286
$ws.Range("I28").Value = 128.86208
$ws.Range("J28").Value = 536.5
$ws.Range("K28").Value = 128.86208
$ws.Range("L28").Value = 536.5
$ws.Range("M28").Value = 356.13792
$ws.Range("N28").Value = -1506.5

$ws.Range("H32").Value = 5808295
$ws.Range("I32").Value = 123.333336
$ws.Range("J32").Value = 7744352
$ws.Range("K32").Value = 123.333336
$ws.Range("L32").Value = 7744352
$ws.Range("M32").Value = 202.666664
$ws.Range("N32").Value = -7745004

$ws.Range("H40").Value = 3866.2354
$ws.Range("I40").Value = 922.2222
$ws.Range("J40").Value = 7178.25
$ws.Range("K40").Value = 922.2222
$ws.Range("L40").Value = 7178.25
$ws.Range("M40").Value = -747.2222
$ws.Range("N40").Value = -7528.25

$ws.Range("H98").Value = 2037.0834
$ws.Range("I98").Value = 1944.1111
$ws.Range("J98").Value = 2316
$ws.Range("K98").Value = 1944.1111
$ws.Range("L98").Value = 2316
$ws.Range("M98").Value = -446.1111000000001
$ws.Range("N98").Value = -5312

$ws.Range("H116").Value = 76644.55499999999
$ws.Range("I116").Value = 114453.266
$ws.Range("K116").Value = 114453.266
$ws.Range("M116").Value = -111011.266

$ws.Range("H122").Value = 2037.0834
$ws.Range("I122").Value = 1944.1111
$ws.Range("J122").Value = 2316
$ws.Range("K122").Value = 5832.3333
$ws.Range("L122").Value = 6948
$ws.Range("M122").Value = -3382.3333
$ws.Range("N122").Value = -11848

$ws.Range("H138").Value = 2721
$ws.Range("J138").Value = 2832.1206
$ws.Range("L138").Value = 8496.361800000001
$ws.Range("N138").Value = -18776.3618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 14346
$ws.Range("I37").Value = 3000
$ws.Range("J37").Value = 20019
$ws.Range("K37").Value = 3000
$ws.Range("L37").Value = 20019
$ws.Range("M37").Value = -2727
$ws.Range("N37").Value = -20565

$ws.Range("H44").Value = 15309.8
$ws.Range("J44").Value = 15309.8
$ws.Range("L44").Value = 15309.8
$ws.Range("N44").Value = -16285.8

$ws.Range("H55").Value = 28000
$ws.Range("J55").Value = 28000
$ws.Range("L55").Value = 28000
$ws.Range("N55").Value = -28630

$ws.Range("H80").Value = 27819
$ws.Range("J80").Value = 27819
$ws.Range("L80").Value = 27819
$ws.Range("N80").Value = -29815

$ws.Range("H83").Value = 27819
$ws.Range("J83").Value = 27819
$ws.Range("L83").Value = 83457
$ws.Range("N83").Value = -93441

$ws.Range("H109").Value = 15800
$ws.Range("J109").Value = 15800
$ws.Range("L109").Value = 15800
$ws.Range("N109").Value = -18574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3451.9565
$ws.Range("I31").Value = 2817.56
$ws.Range("J31").Value = 4207.1904
$ws.Range("K31").Value = 2817.56
$ws.Range("L31").Value = 4207.1904
$ws.Range("M31").Value = -2522.56
$ws.Range("N31").Value = -4797.1904

$ws.Range("H34").Value = 3451.9565
$ws.Range("I34").Value = 2817.56
$ws.Range("J34").Value = 4207.1904
$ws.Range("K34").Value = 2817.56
$ws.Range("L34").Value = 4207.1904
$ws.Range("M34").Value = -2615.56
$ws.Range("N34").Value = -4611.1904

$ws.Range("H99").Value = 86090.164
$ws.Range("I99").Value = 145458.86
$ws.Range("K99").Value = 145458.86
$ws.Range("M99").Value = -143960.86

$ws.Range("H122").Value = 3797.75
$ws.Range("I122").Value = 5870.5
$ws.Range("J122").Value = 1725
$ws.Range("K122").Value = 17611.5
$ws.Range("L122").Value = 5175
$ws.Range("M122").Value = -15161.5
$ws.Range("N122").Value = -10075

$ws.Range("H126").Value = 86090.164
$ws.Range("I126").Value = 145458.86
$ws.Range("K126").Value = 436376.58
$ws.Range("M126").Value = -433906.58

$ws.Range("H132").Value = 11113309
$ws.Range("I132").Value = 14287066
$ws.Range("J132").Value = 5159.8
$ws.Range("K132").Value = 42861198
$ws.Range("L132").Value = 15479.4
$ws.Range("M132").Value = -42858668
$ws.Range("N132").Value = -20539.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1851.4777
$ws.Range("J131").Value = 1500.2034
$ws.Range("L131").Value = 4500.6102
$ws.Range("N131").Value = -14580.6102

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5961.75
$ws.Range("I113").Value = 7406
$ws.Range("J113").Value = 1629
$ws.Range("K113").Value = 7406
$ws.Range("L113").Value = 1629
$ws.Range("M113").Value = -5236
$ws.Range("N113").Value = -5969

$ws.Range("H122").Value = 2101.1428
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 2127
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 6381
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -11281

$ws.Range("H126").Value = 2488.6191
$ws.Range("I126").Value = 2451.2307
$ws.Range("J126").Value = 2549.375
$ws.Range("K126").Value = 7353.6921
$ws.Range("L126").Value = 7648.125
$ws.Range("M126").Value = -4883.6921
$ws.Range("N126").Value = -12588.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2512.0908
$ws.Range("I7").Value = 2993.75
$ws.Range("J7").Value = 2236.8572
$ws.Range("K7").Value = 2993.75
$ws.Range("L7").Value = 2236.8572
$ws.Range("M7").Value = -2881.75
$ws.Range("N7").Value = -2460.8572

$ws.Range("H126").Value = 2512.0908
$ws.Range("I126").Value = 2993.75
$ws.Range("J126").Value = 2236.8572
$ws.Range("K126").Value = 8981.25
$ws.Range("L126").Value = 6710.571599999999
$ws.Range("M126").Value = -6511.25
$ws.Range("N126").Value = -11650.5716

$ws.Range("H128").Value = 36392
$ws.Range("J128").Value = 36392
$ws.Range("L128").Value = 36392
$ws.Range("N128").Value = -46352

$ws.Range("H132").Value = 6542291
$ws.Range("I132").Value = 2251.2354
$ws.Range("K132").Value = 6753.706200000001
$ws.Range("M132").Value = -4223.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 57693
$ws.Range("I122").Value = 64604.625
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 193813.875
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -191363.875
$ws.Range("N122").Value = -12100

$ws.Range("H132").Value = 1640.1409
$ws.Range("I132").Value = 764.7646999999999
$ws.Range("J132").Value = 3872.35
$ws.Range("K132").Value = 2294.2941
$ws.Range("L132").Value = 11617.05
$ws.Range("M132").Value = 235.7058999999999
$ws.Range("N132").Value = -16677.05

$ws.Range("H136").Value = 2027.2222
$ws.Range("I136").Value = 1102.7778
$ws.Range("J136").Value = 4800.5557
$ws.Range("K136").Value = 3308.3334
$ws.Range("L136").Value = 14401.6671
$ws.Range("M136").Value = -758.3334000000004
$ws.Range("N136").Value = -19501.6671
